$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Row 57 becomes a filled "Coding" entry (mirrors rows 53/54/56 formatting).
$ws.Range("A57").Value = "Coding"
$ws.Range("C57").Value = 20
$ws.Range("D57").Value = "Changed the constructors of multiple classes so that Object creation is easier"
$ws.Range("E57").Value = "Finished"

# F57 stores a time-of-day value (13:55) using the same time number format as the other rows.
$ws.Range("F57").Value = 0.57986111111111105
$ws.Range("F57").NumberFormat = $ws.Range("F56").NumberFormat

# Row height grows to match the other "Coding" rows (30pt custom height).
$ws.Rows.Item(57).RowHeight = 30

# Keep the active selection where the author left it after editing row 57.
$ws.Range("D57").Select()
